# Weekly update: a new pair of price records (date 2022-08-19, serial 44792)
# is inserted at the top of the data block (rows 60-61), pushing all the
# existing records (rows 60-123) down by two rows (to rows 62-125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the existing data block (rows 60-123, columns A-R) using
#    Value2 so dates come back as raw numeric serials (not locale strings),
#    which keeps the later bulk write from re-triggering Excel's
#    "looks like a date, pick a default date format" auto-detection.
$srcRange = $ws.Range("A60:R123")
$data = $srcRange.Value2()

# 2) Push that whole block down two rows: old row N (60..123) -> new row N+2 (62..125).
$destRange = $ws.Range("A62:R125")
$destRange.Value = $data

# Make sure the date column keeps its existing date format on the two
# brand new rows at the bottom (124-125) that did not exist before.
$ws.Range("D124:D125").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 3) Write the two brand new records into the now-vacated rows 60-61.
#    (Mercado/Región/Codreg/Categoría/Variedad/Origen/Clasificación columns
#    are constant for this whole sheet.)
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44792
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = 100112036
$ws.Range("G60").Value = "Caigua"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 140
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = 8500
$ws.Range("N60").Value = "$/caja 20 kilos"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 425
$ws.Range("Q60").Value = 20
$ws.Range("R60").Value = "Hortaliza"

$ws.Range("A61").Value = 1
$ws.Range("B61").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C61").Value = "Arica y Parinacota"
$ws.Range("D61").Value = 44792
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = 100112036
$ws.Range("G61").Value = "Caigua"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Segunda"
$ws.Range("J61").Value = 130
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = 6500
$ws.Range("N61").Value = "$/caja 20 kilos"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 325
$ws.Range("Q61").Value = 20
$ws.Range("R61").Value = "Hortaliza"
